# Fixed bug allowing students First and Last name to be entered in SSC
#
# A new student record (Lname="Lyn", Fname="Omari", DOB="11/07",
# DOB old = 11/7/1997) is inserted as the new row 2 of the Table1
# listobject on Sheet1, pushing all the existing data rows down by one
# (table grows from A1:F17 to A1:F18).

$wb  = $excel.ActiveWorkbook
$ws  = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Insert a blank row right below the header (row 2), shifting the
# existing data down.
$ws.Rows.Item(2).Insert()

# Populate the new record.
$ws.Range("A2").Value = "Lyn"
$ws.Range("B2").Value = "Omari"
$ws.Range("C2").Value = "11/07"
$ws.Range("D2").Value = 35741

# The freshly inserted row has no formatting yet - pull it from the row
# right below (the original row 2, now row 3) so the DOB/"DOB old"
# columns keep their existing number formats.
$ws.Range("C3:D3").Copy()
$ws.Range("C2:D2").PasteSpecial(-4122)

# Grow the table so it covers the newly added row.
$tbl.Resize($ws.Range("A1:F18"))

# Match the saved selection state.
$ws.Range("E6").Select() | Out-Null
